# Rewrite the " m:'doc.html'.fromHTMLURI() " M2Doc field (fldChar begin /
# instrText... / fldChar end) into plain literal text runs:
#   { m : ' doc.html <bookmark _GoBack> '.fromHTMLURI() }
# i.e. drop the field-code wrapper and keep the text as ordinary w:t runs,
# split the same way the original instrText runs were split, and keep the
# _GoBack bookmark sitting between "doc.html" and "'.fromHTMLURI()".

$d = $word.ActiveDocument

# Find the paragraph that holds the lone field (the query paragraph).
$fieldPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $fieldPara = $p
    }
}

$field = $fieldPara.Range.Fields.Item(1)

# WordprocessingML fragment that replaces the paragraph's content. Each
# token gets its own <w:r><w:t>...</w:t></w:r> (mirroring the original
# instrText run split), the _GoBack bookmark stays between "doc.html" and
# "'.fromHTMLURI()", and the closing "}" keeps xml:space="preserve" (it
# used to be the trailing-space instrText run right before fldChar end).
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$q = "'"
$xmlFrag = '<w:p ' + $ns + `
    ' w:rsidR="' + $fieldPara.Range.ParagraphFormat.Parent.Information(3) + '"' + `
    '><w:r><w:t>{</w:t></w:r>' + `
    '<w:r><w:t>m</w:t></w:r>' + `
    '<w:r><w:t>:</w:t></w:r>' + `
    '<w:r><w:t>' + $q + '</w:t></w:r>' + `
    '<w:r><w:t>doc.html</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t>' + $q + '.fromHTMLURI()</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">}</w:t></w:r></w:p>'

$fieldPara.Range.InsertXML($xmlFrag)
